$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Populations")

$ws.Range("A2").Value = "POP_EU_ADULT_2022"
$ws.Range("B2").Value = "Adults EU 2022"
$ws.Range("C2").Value = "EU adult population (synthetic) for 2022"
$ws.Range("D2").Value = "EU"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("E2").Value = "2022-01-01"
$ws.Range("F2").Value = "2022-12-31"
$ws.Range("H2").Value = 1000000
